$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LabelsUsedMultipleTerms")
Write-Host $ws.Name
